$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 978.625
$ws.Range("J17").Value = 977.2
$ws.Range("L17").Value = 2931.6
$ws.Range("N17").Value = -3267.6
$ws.Range("H40").Value = 3350
$ws.Range("I40").Value = 1866.6666
$ws.Range("J40").Value = 4833.3335
$ws.Range("K40").Value = 1866.6666
$ws.Range("L40").Value = 4833.3335
$ws.Range("M40").Value = -1691.6666
$ws.Range("N40").Value = -5183.3335
$ws.Range("H86").Value = 2379.0588
$ws.Range("I86").Value = 1944.6666
$ws.Range("K86").Value = 1944.6666
$ws.Range("M86").Value = -821.6666
$ws.Range("H89").Value = 2379.0588
$ws.Range("I89").Value = 1944.6666
$ws.Range("K89").Value = 9723.333000000001
$ws.Range("M89").Value = -4107.333000000001
$ws.Range("H96").Value = 5495357.5
$ws.Range("I96").Value = 11905020
$ws.Range("J96").Value = 1361.2858
$ws.Range("K96").Value = 35715060
$ws.Range("L96").Value = 4083.8574
$ws.Range("M96").Value = -35713687
$ws.Range("N96").Value = -6829.857400000001
$ws.Range("H98").Value = 3155.0715
$ws.Range("I98").Value = 1858.8462
$ws.Range("K98").Value = 1858.8462
$ws.Range("M98").Value = -360.8462
$ws.Range("H116").Value = 21622.277
$ws.Range("I116").Value = 24771.643
$ws.Range("J116").Value = 10599.5
$ws.Range("K116").Value = 24771.643
$ws.Range("L116").Value = 10599.5
$ws.Range("M116").Value = -21329.643
$ws.Range("N116").Value = -17483.5
$ws.Range("H122").Value = 3155.0715
$ws.Range("I122").Value = 1858.8462
$ws.Range("K122").Value = 5576.5386
$ws.Range("M122").Value = -3126.5386
$ws.Range("H132").Value = 22244.354
$ws.Range("I132").Value = 23942.193
$ws.Range("J132").Value = 4700
$ws.Range("K132").Value = 71826.579
$ws.Range("L132").Value = 14100
$ws.Range("M132").Value = -69296.579
$ws.Range("N132").Value = -19160
$ws.Range("H137").Value = 43833.57
$ws.Range("I137").Value = 31126.5
$ws.Range("J137").Value = 75601.25
$ws.Range("K137").Value = 93379.5
$ws.Range("L137").Value = 226803.75
$ws.Range("M137").Value = -90829.5
$ws.Range("N137").Value = -231903.75
$ws.Range("H138").Value = 28863.816
$ws.Range("I138").Value = 1989.5416
$ws.Range("K138").Value = 5968.6248
$ws.Range("M138").Value = -828.6247999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24529.977
$ws.Range("I32").Value = 27498.05
$ws.Range("K32").Value = 27498.05
$ws.Range("M32").Value = -27211.05
$ws.Range("H61").Value = 8127.6
$ws.Range("I61").Value = 1325
$ws.Range("J61").Value = 35338
$ws.Range("K61").Value = 1325
$ws.Range("L61").Value = 35338
$ws.Range("M61").Value = -1113
$ws.Range("N61").Value = -35762
$ws.Range("H63").Value = 2893.4167
$ws.Range("I63").Value = 1965.375
$ws.Range("K63").Value = 1965.375
$ws.Range("M63").Value = -1279.375
$ws.Range("H66").Value = 2893.4167
$ws.Range("I66").Value = 1965.375
$ws.Range("K66").Value = 9826.875
$ws.Range("M66").Value = -6394.875
$ws.Range("H74").Value = 872805.1
$ws.Range("I74").Value = 1500705.5
$ws.Range("J74").Value = 35604.668
$ws.Range("K74").Value = 1500705.5
$ws.Range("L74").Value = 35604.668
$ws.Range("M74").Value = -1499831.5
$ws.Range("N74").Value = -37352.668
$ws.Range("H77").Value = 872805.1
$ws.Range("I77").Value = 1500705.5
$ws.Range("J77").Value = 35604.668
$ws.Range("K77").Value = 7503527.5
$ws.Range("L77").Value = 178023.34
$ws.Range("M77").Value = -7499159.5
$ws.Range("N77").Value = -186759.34
$ws.Range("H102").Value = 1981.3043
$ws.Range("I102").Value = 1552
$ws.Range("K102").Value = 1552
$ws.Range("M102").Value = 70
$ws.Range("H136").Value = 8127.6
$ws.Range("I136").Value = 1325
$ws.Range("J136").Value = 35338
$ws.Range("K136").Value = 3975
$ws.Range("L136").Value = 106014
$ws.Range("M136").Value = -1425
$ws.Range("N136").Value = -111114

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2081.6316
$ws.Range("J105").Value = 2825.25
$ws.Range("L105").Value = 2825.25
$ws.Range("N105").Value = -6319.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 61666.668
$ws.Range("I69").Value = 52500
$ws.Range("J69").Value = 80000
$ws.Range("K69").Value = 52500
$ws.Range("L69").Value = 80000
$ws.Range("M69").Value = -51751
$ws.Range("N69").Value = -81498
$ws.Range("H72").Value = 61666.668
$ws.Range("I72").Value = 52500
$ws.Range("J72").Value = 80000
$ws.Range("K72").Value = 157500
$ws.Range("L72").Value = 240000
$ws.Range("M72").Value = -153756
$ws.Range("N72").Value = -247488
$ws.Range("H86").Value = 31814.518
$ws.Range("I86").Value = 44360.234
$ws.Range("J86").Value = 14041.417
$ws.Range("K86").Value = 44360.234
$ws.Range("L86").Value = 14041.417
$ws.Range("M86").Value = -43237.234
$ws.Range("N86").Value = -16287.417
$ws.Range("H89").Value = 31814.518
$ws.Range("I89").Value = 44360.234
$ws.Range("J89").Value = 14041.417
$ws.Range("K89").Value = 221801.17
$ws.Range("L89").Value = 70207.08499999999
$ws.Range("M89").Value = -216185.17
$ws.Range("N89").Value = -81439.08499999999
$ws.Range("H132").Value = 57022.832
$ws.Range("I132").Value = 67560.734
$ws.Range("J132").Value = 4333.3335
$ws.Range("K132").Value = 202682.202
$ws.Range("L132").Value = 13000.0005
$ws.Range("M132").Value = -200152.202
$ws.Range("N132").Value = -18060.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 821
$ws.Range("I2").Value = 821
$ws.Range("K2").Value = 4926
$ws.Range("M2").Value = -4813
$ws.Range("H3").Value = 1022.4286
$ws.Range("I3").Value = 1022.4286
$ws.Range("K3").Value = 3067.2858
$ws.Range("M3").Value = -2955.2858
$ws.Range("H12").Value = 189
$ws.Range("J12").Value = 196.21428
$ws.Range("L12").Value = 588.64284
$ws.Range("N12").Value = -934.64284
$ws.Range("H117").Value = 1416
$ws.Range("J117").Value = 1416
$ws.Range("L117").Value = 4248
$ws.Range("N117").Value = -11132
$ws.Range("H131").Value = 114302.19
$ws.Range("I131").Value = 429105.1
$ws.Range("J131").Value = 1872.5714
$ws.Range("K131").Value = 1287315.3
$ws.Range("L131").Value = 5617.7142
$ws.Range("M131").Value = -1282275.3
$ws.Range("N131").Value = -15697.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 74282.5
$ws.Range("J95").Value = 74282.5
$ws.Range("L95").Value = 74282.5
$ws.Range("N95").Value = -79774.5
$ws.Range("H97").Value = 1012.1724
$ws.Range("I97").Value = 898.2308
$ws.Range("K97").Value = 898.2308
$ws.Range("M97").Value = -402.2308

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1984.9524
$ws.Range("I16").Value = 1275.8823
$ws.Range("J16").Value = 4998.5
$ws.Range("K16").Value = 1275.8823
$ws.Range("L16").Value = 4998.5
$ws.Range("M16").Value = -1105.8823
$ws.Range("N16").Value = -5338.5
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H99").Value = 42666.332
$ws.Range("I99").Value = 42666.332
$ws.Range("K99").Value = 42666.332
$ws.Range("M99").Value = -39671.332
$ws.Range("H100").Value = 2730.6316
$ws.Range("I100").Value = 2506.3076
$ws.Range("J100").Value = 3216.6667
$ws.Range("K100").Value = 2506.3076
$ws.Range("L100").Value = 3216.6667
$ws.Range("M100").Value = -1965.3076
$ws.Range("N100").Value = -4298.6667
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7748.4614
$ws.Range("I81").Value = 8144.1665
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 16288.333
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -15227.333
$ws.Range("N81").Value = -8122
$ws.Range("H84").Value = 7748.4614
$ws.Range("I84").Value = 8144.1665
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 81441.66500000001
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -76137.66500000001
$ws.Range("N84").Value = -40608
$ws.Range("H122").Value = 36164.824
$ws.Range("I122").Value = 40861.715
$ws.Range("J122").Value = 3286.6
$ws.Range("K122").Value = 122585.145
$ws.Range("L122").Value = 9859.799999999999
$ws.Range("M122").Value = -120135.145
$ws.Range("N122").Value = -14759.8
$ws.Range("H132").Value = 39219.79
$ws.Range("I132").Value = 51237.57
$ws.Range("K132").Value = 153712.71
$ws.Range("M132").Value = -151182.71
